$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Move the old "expectedHeader" / "Where To Buy" pair into the new column D,
# and introduce a new column C carrying "expectedUrl" / "where-to-buy".
$ws.Range("D1").Value = "expectedHeader"
$ws.Range("D2").Value = "Where To Buy"

$ws.Range("C1").Value = "expectedUrl"
$ws.Range("C2").Value = "where-to-buy"

$ws.Range("A2").Value = "Verify 'Where To Buy' header"

$ws.Range("A2").Select()
